$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 268, pushing existing rows 268-272 down to 270-274
$ws.Range("A268:A269").EntireRow.Insert()

# Fill new row 268 - Carson
$ws.Cells.Item(268, 1).Value = 5
$ws.Cells.Item(268, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(268, 3).Value = "Maule"
$ws.Cells.Item(268, 4).Value = 44578
$ws.Cells.Item(268, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(268, 5).Value = 7
$ws.Cells.Item(268, 6).Value = "Fruta"
$ws.Cells.Item(268, 7).Value = 100103
$ws.Cells.Item(268, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(268, 9).Value = 100103004
$ws.Cells.Item(268, 10).Value = "Durazno"
$ws.Cells.Item(268, 11).Value = "Carson"
$ws.Cells.Item(268, 12).Value = "Especial"
$ws.Cells.Item(268, 13).Value = 300
$ws.Cells.Item(268, 14).Value = 13000
$ws.Cells.Item(268, 15).Value = 13000
$ws.Cells.Item(268, 16).Value = 13000
$ws.Cells.Item(268, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(268, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(268, 19).Value = 867
$ws.Cells.Item(268, 20).Value = 15

# Fill new row 269 - Elegant Lady
$ws.Cells.Item(269, 1).Value = 5
$ws.Cells.Item(269, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(269, 3).Value = "Maule"
$ws.Cells.Item(269, 4).Value = 44578
$ws.Cells.Item(269, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(269, 5).Value = 7
$ws.Cells.Item(269, 6).Value = "Fruta"
$ws.Cells.Item(269, 7).Value = 100103
$ws.Cells.Item(269, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(269, 9).Value = 100103004
$ws.Cells.Item(269, 10).Value = "Durazno"
$ws.Cells.Item(269, 11).Value = "Elegant Lady"
$ws.Cells.Item(269, 12).Value = "Especial"
$ws.Cells.Item(269, 13).Value = 180
$ws.Cells.Item(269, 14).Value = 14000
$ws.Cells.Item(269, 15).Value = 14000
$ws.Cells.Item(269, 16).Value = 14000
$ws.Cells.Item(269, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(269, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(269, 19).Value = 933
$ws.Cells.Item(269, 20).Value = 15
